$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (12S / DADA2Spec / 50%)
$ws.Range("D11").Value = 30
$ws.Range("E11").Value = 3.5
$ws.Range("G11").Value = 65.5
$ws.Range("H11").Value = 0.31
$ws.Range("I11").Value = 0.9
$ws.Range("J11").Value = 0.47
$ws.Range("K11").Value = 0.65
$ws.Range("L11").Value = 0.3

# Row 12 (12S / DADA2Spec / 30%)
$ws.Range("D12").Value = 40.1
$ws.Range("E12").Value = 4.1
$ws.Range("G12").Value = 54.8
$ws.Range("H12").Value = 0.42
$ws.Range("I12").Value = 0.91
$ws.Range("J12").Value = 0.58
$ws.Range("K12").Value = 0.74
$ws.Range("L12").Value = 0.41

# Row 13 (12S / DADA2Spec / 70%)
$ws.Range("D13").Value = 16.3
$ws.Range("E13").Value = 1.5
$ws.Range("G13").Value = 81.2
$ws.Range("H13").Value = 0.17
$ws.Range("I13").Value = 0.92
$ws.Range("J13").Value = 0.28
$ws.Range("K13").Value = 0.48
$ws.Range("L13").Value = 0.16

# Row 53 (16S / DADA2Spec / 50%)
$ws.Range("D53").Value = 24.5
$ws.Range("E53").Value = 0.4
$ws.Range("G53").Value = 74.09999999999999
$ws.Range("H53").Value = 0.25
$ws.Range("I53").Value = 0.98
$ws.Range("J53").Value = 0.4
$ws.Range("K53").Value = 0.62
$ws.Range("L53").Value = 0.25

# Row 54 (16S / DADA2Spec / 30%)
$ws.Range("D54").Value = 36.2
$ws.Range("E54").Value = 0.9
$ws.Range("G54").Value = 61.9
$ws.Range("H54").Value = 0.37
$ws.Range("I54").Value = 0.98
$ws.Range("J54").Value = 0.54
$ws.Range("K54").Value = 0.73
$ws.Range("L54").Value = 0.37

# Row 55 (16S / DADA2Spec / 70%)
$ws.Range("D55").Value = 14.7
$ws.Range("E55").Value = 0.2
$ws.Range("G55").Value = 84.09999999999999
$ws.Range("H55").Value = 0.15
$ws.Range("I55").Value = 0.99
$ws.Range("J55").Value = 0.26
$ws.Range("K55").Value = 0.46
$ws.Range("L55").Value = 0.15
